$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A38").Value = 'BleepingComputer'
$ws.Range("B38").Value = 'Ivanti Avalanche impacted by critical pre-auth stack buffer overflows'
$ws.Range("C38").Value = 'https://www.bleepingcomputer.com/news/security/ivanti-avalanche-impacted-by-critical-pre-auth-stack-buffer-overflows/'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2023-08-15'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = 'Two stack-based buffer overflows collectively tracked as CVE-2023-32560 impact Ivanti Avalanche, an enterprise mobility management (EMM) solution designed to manage, monitor, and secure a wide range of mobile devices. [...]'

$ws.Range("A39").Value = 'BleepingComputer'
$ws.Range("B39").Value = 'LinkedIn accounts hacked in widespread hijacking campaign'
$ws.Range("C39").Value = 'https://www.bleepingcomputer.com/news/security/linkedin-accounts-hacked-in-widespread-hijacking-campaign/'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2023-08-15'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = 'LinkedIn is being targeted in a wave of account hacks resulting in many accounts being locked out for security reasons or ultimately hijacked by attackers. [...]'

$ws.Range("A40").Value = 'BleepingComputer'
$ws.Range("B40").Value = 'Almost 2,000 Citrix NetScaler servers backdoored in hacking campaign'
$ws.Range("C40").Value = 'https://www.bleepingcomputer.com/news/security/almost-2-000-citrix-netscaler-servers-backdoored-in-hacking-campaign/'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2023-08-15'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = 'A threat actor has compromised close to 2,000 thousand Citrix NetScaler servers in a massive campaign exploiting the critical-severity remote code execution tracked as CVE-2023-3519. [...]'

$ws.Range("A41").Value = 'BleepingComputer'
$ws.Range("B41").Value = 'Raccoon Stealer malware returns with new stealthier version'
$ws.Range("C41").Value = 'https://www.bleepingcomputer.com/news/security/raccoon-stealer-malware-returns-with-new-stealthier-version/'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2023-08-15'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = 'The developers of Raccoon Stealer information-stealing malware have ended their 6-month hiatus from hacker forums to promote a new 2.3.0 version of the malware to cyber criminals. [...]'

$ws.Range("A42").Value = 'BleepingComputer'
$ws.Range("B42").Value = 'New CVE-2023-3519 scanner detects hacked Citrix ADC, Gateway devices'
$ws.Range("C42").Value = 'https://www.bleepingcomputer.com/news/security/new-cve-2023-3519-scanner-detects-hacked-citrix-adc-gateway-devices/'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2023-08-15'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = 'Mandiant has released a scanner to check if a Citrix NetScaler Application Delivery Controller (ADC) or NetScaler Gateway Appliance was compromised in widespread attacks exploiting the CVE-2023-3519 vulnerability.  [...]'

$ws.Range("A43").Value = 'BleepingComputer'
$ws.Range("B43").Value = 'Back to school security against ransomware attacks on K-12 and colleges'
$ws.Range("C43").Value = 'https://www.bleepingcomputer.com/news/security/back-to-school-security-against-ransomware-attacks-on-k-12-and-colleges/'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2023-08-15'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = 'As we get back to school, K-12 and colleges are increasingly at risk from ransomware and data theft attacks. Learn more from Specops Software on the steps IT teams at education institutes can take to protect their care orgs from disruption and stolen data. [...]'

$ws.Range("A44").Value = 'BleepingComputer'
$ws.Range("B44").Value = 'Threat actors use beta apps to bypass mobile app store security'
$ws.Range("C44").Value = 'https://www.bleepingcomputer.com/news/security/threat-actors-use-beta-apps-to-bypass-mobile-app-store-security/'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2023-08-14'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = 'The FBI is warning of a new tactic used by cybercriminals where they promote malicious beta versions of cryptocurrency investment apps on popular mobile app stores that are then used to steal crypto. [...]'

$ws.Range("A45").Value = 'BleepingComputer'
$ws.Range("B45").Value = 'Discord.io confirms breach after hacker steals data of 760K users'
$ws.Range("C45").Value = 'https://www.bleepingcomputer.com/news/security/discordio-confirms-breach-after-hacker-steals-data-of-760k-users/'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2023-08-14'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = 'The Discord.io custom invite service has temporarily shut down after suffering a data breach exposing the information of 760,000 members. [...]'

$ws.Range("A46").Value = 'BleepingComputer'
$ws.Range("B46").Value = 'Over 100K hacking forums accounts exposed by info-stealing malware'
$ws.Range("C46").Value = 'https://www.bleepingcomputer.com/news/security/over-100k-hacking-forums-accounts-exposed-by-info-stealing-malware/'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2023-08-14'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = 'Researchers discovered 120,000 infected systems that contained credentials for cybercrime forums. Many of the computers belong to hackers, the researchers say. [...]'

$ws.Range("A47").Value = 'BleepingComputer'
$ws.Range("B47").Value = 'Microsoft enables Windows Kernel CVE-2023-32019 fix for everyone'
$ws.Range("C47").Value = 'https://www.bleepingcomputer.com/news/microsoft/microsoft-enables-windows-kernel-cve-2023-32019-fix-for-everyone/'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2023-08-14'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = 'Microsoft has enabled a fix for a Kernel information disclosure vulnerability by default for everyone after previously disabling it out of concerns it could introduce breaking changes to Windows. [...]'

$ws.Range("A48").Value = 'BleepingComputer'
$ws.Range("B48").Value = 'FBI warns of increasing cryptocurrency recovery scams'
$ws.Range("C48").Value = 'https://www.bleepingcomputer.com/news/security/fbi-warns-of-increasing-cryptocurrency-recovery-scams/'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2023-08-14'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = 'The FBI is warning of an increase in scammers pretending to be recovery companies that can help victims of cryptocurrency investment scams recover lost assets. [...]'

$ws.Range("A49").Value = 'BleepingComputer'
$ws.Range("B49").Value = 'Monti ransomware targets VMware ESXi servers with new Linux locker'
$ws.Range("C49").Value = 'https://www.bleepingcomputer.com/news/security/monti-ransomware-targets-vmware-esxi-servers-with-new-linux-locker/'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2023-08-14'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = 'The Monti ransomware has returned to action after a two-month hiatus, now targeting primarily legal and government organizations, and VMware ESXi servers using a new Linux variant that is vastly different from its predecessors. [...]'

$ws.Range("A50").Value = 'BleepingComputer'
$ws.Range("B50").Value = 'Criminal IP Teams Up with PolySwarm to Strengthen Threat Detection'
$ws.Range("C50").Value = 'https://www.bleepingcomputer.com/news/security/criminal-ip-teams-up-with-polyswarm-to-strengthen-threat-detection/'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2023-08-14'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = 'The addition of Criminal IP as a new contributor to PolySwarm''s malicious URL detection represents a significant leap in specialized threat identification. Learn more from Criminal IP about this new collaboration. [...]'
